$wb = $excel.ActiveWorkbook

# --- Update product descriptions / material codes across sheets ---
# These are the content-level edits behind this commit: several catalog
# items were re-specified (new codes, and the submersible pump line item
# was replaced with a different spec entirely).

$ws1 = $wb.Worksheets.Item("MWSS with Submersible Pump")
$ws1.Range("A3").Value = "Supply of Single Phase Submersible Pump 2 HP, Discharge 5200-10600 lph, Head 57-18m without Panel Board (Code: GWDMR149)"
$ws1.Range("A4").Value = "40 mm dia UPVC pipe (Code: GWDMR073)"
$ws1.Range("A5").Value = "14 mm thick nylon rope (Code: GWDMR069)"
$ws1.Range("A9").Value = "40 mm SS Adapter (Code: GWDMR080)"

$ws2 = $wb.Worksheets.Item("MWSS with Compressor Pump")
$ws2.Range("A3").Value = "Compressor pump 2 HP single phase (Code: GWDMR065)"
$ws2.Range("A10").Value = "14 mm thick nylon rope (Code: GWDMR069)"
$ws2.Range("A14").Value = "40 mm SS Adapter (Code: GWDMR080)"

$ws5 = $wb.Worksheets.Item("150 mm Tubewell Construction")
$ws5.Range("A7").Value = "MS Casing pipe 450 mm dia, 6 mm thickness (Code: GWDMR087)"

$ws6 = $wb.Worksheets.Item("200 mm Tubewell Construction")
$ws6.Range("A7").Value = "MS Casing pipe 450 mm dia, 6 mm thickness (Code: GWDMR087)"

$ws7 = $wb.Worksheets.Item("Submersible Pump Installation")
$ws7.Range("A3").Value = "Supply of Single Phase Submersible Pump 2 HP, Discharge 5200-10600 lph, Head 57-18m without Panel Board (Code: GWDMR149)"
$ws7.Range("A4").Value = "40 mm dia UPVC pipe (Code: GWDMR073)"
$ws7.Range("A5").Value = "14 mm thick nylon rope (Code: GWDMR069)"
$ws7.Range("A9").Value = "40 mm SS Adapter (Code: GWDMR080)"

$ws8 = $wb.Worksheets.Item("Compressor Pump Installation")
$ws8.Range("A6").Value = "20 mm dia HDPE pipe (DG) (8kg) (Code: GWDMR076)"
$ws8.Range("A7").Value = "32 mm dia UPVC pipe (Code: GWDMR072)"
$ws8.Range("A10").Value = "14 mm thick nylon rope (Code: GWDMR069)"
$ws8.Range("A14").Value = "40 mm SS Adapter (Code: GWDMR080)"
